# Update "想去人数" (column F) counts that changed between data refreshes,
# on both the "展览" sheet and the "全部类型" sheet (which aggregates all
# event types and therefore carries the same rows, offset by one due to
# the extra "演出" row inserted at row 14).

$wb = $excel.ActiveWorkbook

$exhibition = $wb.Worksheets.Item("展览")
$exhibition.Range("F2").Value = 11
$exhibition.Range("F3").Value = 1114
$exhibition.Range("F5").Value = 92
$exhibition.Range("F7").Value = 60
$exhibition.Range("F8").Value = 11401
$exhibition.Range("F9").Value = 4321
$exhibition.Range("F11").Value = 31
$exhibition.Range("F13").Value = 2525
$exhibition.Range("F15").Value = 122
$exhibition.Range("F16").Value = 25
$exhibition.Range("F18").Value = 497
$exhibition.Range("F19").Value = 11275
$exhibition.Range("F20").Value = 11153
$exhibition.Range("F25").Value = 39

$allTypes = $wb.Worksheets.Item("全部类型")
$allTypes.Range("F2").Value = 11
$allTypes.Range("F3").Value = 1114
$allTypes.Range("F5").Value = 92
$allTypes.Range("F7").Value = 60
$allTypes.Range("F8").Value = 11401
$allTypes.Range("F9").Value = 4321
$allTypes.Range("F11").Value = 31
$allTypes.Range("F13").Value = 2525
$allTypes.Range("F16").Value = 122
$allTypes.Range("F17").Value = 25
$allTypes.Range("F19").Value = 497
$allTypes.Range("F20").Value = 11275
$allTypes.Range("F21").Value = 11153
$allTypes.Range("F26").Value = 39
